$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the phone number value (and its cell record) that was previously in B2
$ws.Range("B2").Clear()

# Update the active cell selection (cosmetic, matches author's last cursor position)
$ws.Range("B8").Select()
